$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 74107.74448106618
$ws.Range("B3").Value = 70103.87168517186
$ws.Range("B4").Value = 66740.12165126778
$ws.Range("B5").Value = 63950.19071297276
$ws.Range("B6").Value = 62815.86644126129
$ws.Range("B7").Value = 64397.29147207199
$ws.Range("B8").Value = 65143.17165758974
$ws.Range("B9").Value = 72217.35464526717
$ws.Range("B10").Value = 88145.49880963346
$ws.Range("B11").Value = 97351.48440994957
$ws.Range("B12").Value = 101289.2699103402
$ws.Range("B13").Value = 102719.3602059575
$ws.Range("B14").Value = 103954.4797775987
$ws.Range("B15").Value = 107916.4737841847
$ws.Range("B16").Value = 109087.1281679656
$ws.Range("B17").Value = 106697.2795579919
$ws.Range("B18").Value = 100528.1278129718
$ws.Range("B19").Value = 92800.63566401665
$ws.Range("B20").Value = 90054.59211346184
$ws.Range("B21").Value = 86562.54541493443
$ws.Range("B22").Value = 84079.22294763885
$ws.Range("B23").Value = 80850.41620926266
$ws.Range("B24").Value = 78627.59173169297
$ws.Range("B25").Value = 74369.82822427993
